# Add "hydrogen combined cycle" as a new power plant type row, and
# rename the existing "hydrogen" row to "hydrogen combustion turbine"
# on the PTUfIGaMDC sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PTUfIGaMDC")

# Rename the shared string "hydrogen" -> "hydrogen combustion turbine"
# (updates the existing shared-string entry used by A24 in place).
$ws.Range("A24").Value = "hydrogen combustion turbine"

# Add the new row 25: "hydrogen combined cycle" with a 0 flag.
$ws.Range("A25").Value = "hydrogen combined cycle"
$ws.Range("B25").Value = 0

# Give A24 its new formatting (black font color + vertically centered),
# matching the new cellXfs/font introduced for this edit, then copy that
# exact formatting onto A25 so both rows share a single new style entry.
$a24 = $ws.Range("A24")
$a24.Font.Color = 0
$a24.VerticalAlignment = -4108
$a24.Copy()
$ws.Range("A25").PasteSpecial(-4122)

# Update the sheet's selection to B25, then restore the originally
# active sheet ("About") so the workbook's active-tab state is unchanged.
$ws.Activate()
$ws.Range("B25").Select()
$wb.Worksheets.Item("About").Activate()
